$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header revision/date text updates ---
$ws.Range("I2").Value = "Revisão: 1.2"
$ws.Range("I3").Value = "Data: 07/02/2021"

# --- Total % legend on header row ---
$ws.Range("N6").Value = " (%) Total (600%)"

# --- Capture the existing "Testes" task (currently row 13) before overwriting it ---
$oldC = $ws.Range("C13").Value2
$oldD = $ws.Range("D13").Value2
$oldE = $ws.Range("E13").Value2
$oldG = $ws.Range("G13").Value2
$oldH = $ws.Range("H13").Value2

# Give row 14 (today an empty placeholder task slot) the same look as row 13 currently has,
# then move the "Testes" task data down into it.
$ws.Range("B13:O13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("C14").Value = $oldC
$ws.Range("D14").Value = $oldD
$ws.Range("E14").Value = $oldE
$ws.Range("G14").Value = $oldG
$ws.Range("H14").Value = $oldH

# --- Row 13 becomes the new task: "Desenvolver ducumentação" (ADM / Murilo) ---
# Give it the same look as the row above it (row 12 - a normal task row).
$ws.Range("C12:O12").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("C13").Value = "Desenvolver ducumentação"
$ws.Range("D13").Value = "ADM"
$ws.Range("E13").Value = "Murilo"
$ws.Range("G13").Value = 44228
$ws.Range("H13").Value = 44246
$ws.Range("N13").Value = 0

# Column F narrows now that E/F are no longer a single wide merged-width pair.
$ws.Columns("F").ColumnWidth = 5.75

# Match the author's final selection.
$ws.Range("E18").Select()
